# "Generate Report for Handback" - populate Latest Target File / Latest Handback File
# columns (F/G) for the zh-cn and de-de handback reports, update the handback
# timestamps, and flip the Overview/Status text from "Ready for handoff" to
# "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# ---------------------------------------------------------------------------
# 1. Status text now reflects a completed handback instead of a pending
#    handoff. This shared string is used on all three sheets.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$ws1.Range("B2").Value2 = $newStatus
$ws1.Range("C2").Value2 = $newStatus
$ws1.Range("B3").Value2 = $newStatus
$ws1.Range("C3").Value2 = $newStatus
$ws2.Range("C2").Value2 = $newStatus
$ws2.Range("C3").Value2 = $newStatus
$ws3.Range("C2").Value2 = $newStatus
$ws3.Range("C3").Value2 = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn (sheet 2): fill in "Latest Target File" (F) and "Latest Handback
#    File" (G) for both data rows, and stamp the real handback datetime into
#    "Latest Handback DateTime" (H), replacing the 0001-01-01 placeholder.
# ---------------------------------------------------------------------------
$zhXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$ws2.Range("F2").Value2 = "a.md"
$ws2.Range("G2").Value2 = $zhXlf
$ws2.Range("F3").Value2 = "a.md"
$ws2.Range("G3").Value2 = $zhXlf

$ws2.Range("H2").Value2 = "2016-03-20 06:26:07"
$ws2.Range("H3").Value2 = "2016-03-20 06:26:07"

# ---------------------------------------------------------------------------
# 3. de-de (sheet 3): same treatment, with its own handback datetime.
# ---------------------------------------------------------------------------
$deXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$ws3.Range("F2").Value2 = "a.md"
$ws3.Range("G2").Value2 = $deXlf
$ws3.Range("F3").Value2 = "a.md"
$ws3.Range("G3").Value2 = $deXlf

$ws3.Range("H2").Value2 = "2016-03-20 06:26:13"
$ws3.Range("H3").Value2 = "2016-03-20 06:26:13"

# ---------------------------------------------------------------------------
# 4. Give the new F/G cells the same "HyperLink" look as the other linked
#    file-name columns, and wire up real hyperlinks for them. The existing
#    hyperlinks are recreated afterwards (in column order) so the resulting
#    <hyperlinks> list stays ordered the way a row-by-row report generator
#    would emit it.
# ---------------------------------------------------------------------------
function Add-ReportHyperlink($ws, $cellRef, $url, $display) {
    $ws.Range($cellRef).Style = "HyperLink"
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $display) | Out-Null
}

$srcBase = "https://github.com/OpenLocalizationTest/oltest/blob/88c73f87c148680de8005dc81ecda626471b6230/e2e"
$zhHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/82f8f2229fea7b1449a700eeffaf1ac90473b5fb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/440c34920efa11f31ec68eeadb7273288ece9e30/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

$ws2.Hyperlinks.Delete()
Add-ReportHyperlink $ws2 "A2" "$srcBase/a.md" "a.md"
Add-ReportHyperlink $ws2 "B2" "$srcBase/a.md" ".md"
Add-ReportHyperlink $ws2 "D2" "$zhHandoffBase/$zhXlf" $zhXlf
Add-ReportHyperlink $ws2 "F2" "$srcBase/a.md" "a.md"
Add-ReportHyperlink $ws2 "G2" "$zhHandoffBase/$zhXlf" $zhXlf
Add-ReportHyperlink $ws2 "A3" "$srcBase/b.md" "b.md"
Add-ReportHyperlink $ws2 "B3" "$srcBase/b.md" ".md"
Add-ReportHyperlink $ws2 "D3" "$zhHandoffBase/$zhXlf" $zhXlf
Add-ReportHyperlink $ws2 "F3" "$srcBase/a.md" "a.md"
Add-ReportHyperlink $ws2 "G3" "$zhHandoffBase/$zhXlf" $zhXlf

$ws3.Hyperlinks.Delete()
Add-ReportHyperlink $ws3 "A2" "$srcBase/a.md" "a.md"
Add-ReportHyperlink $ws3 "B2" "$srcBase/a.md" ".md"
Add-ReportHyperlink $ws3 "D2" "$deHandoffBase/$deXlf" $deXlf
Add-ReportHyperlink $ws3 "F2" "$srcBase/a.md" "a.md"
Add-ReportHyperlink $ws3 "G2" "$deHandoffBase/$deXlf" $deXlf
Add-ReportHyperlink $ws3 "A3" "$srcBase/b.md" "b.md"
Add-ReportHyperlink $ws3 "B3" "$srcBase/b.md" ".md"
Add-ReportHyperlink $ws3 "D3" "$deHandoffBase/$deXlf" $deXlf
Add-ReportHyperlink $ws3 "F3" "$srcBase/a.md" "a.md"
Add-ReportHyperlink $ws3 "G3" "$deHandoffBase/$deXlf" $deXlf

Write-Output "Handback report generated."
